$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 49 (Primera, 2022-04-08) with the new weekly prices (2022-07-11)
$ws.Range("D49").Value = 44753
$ws.Range("K49").Value = 600
$ws.Range("L49").Value = 700
$ws.Range("M49").Value = 650
$ws.Range("P49").Value = 650

# Insert a new row at 50 (pushes the old rows 50 down to 51) to hold the new
# "Segunda" quality entry for the same 2022-07-11 week
$ws.Rows.Item(50).Insert()

$ws.Range("A50").Value = 7
$ws.Range("B50").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C50").Value = "Ñuble"
$ws.Range("D50").Value = 44753
$ws.Range("E50").Value = 16
$ws.Range("F50").Value = 100112040
$ws.Range("G50").Value = "Cilantro"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Segunda"
$ws.Range("J50").Value = 100
$ws.Range("K50").Value = 500
$ws.Range("L50").Value = 500
$ws.Range("M50").Value = 500
$ws.Range("N50").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O50").Value = "Provincia de Diguillín"
$ws.Range("P50").Value = 500
$ws.Range("Q50").Value = 1
$ws.Range("R50").Value = "Hortaliza"

# Insert another row at 51 (pushes the former row 51 [old row 50] down to 52)
# to hold the prior "Primera" 2022-04-08 reading that's being preserved as
# its own historical row
$ws.Rows.Item(51).Insert()

$ws.Range("A51").Value = 7
$ws.Range("B51").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C51").Value = "Ñuble"
$ws.Range("D51").Value = 44659
$ws.Range("E51").Value = 16
$ws.Range("F51").Value = 100112040
$ws.Range("G51").Value = "Cilantro"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 200
$ws.Range("K51").Value = 550
$ws.Range("L51").Value = 600
$ws.Range("M51").Value = 575
$ws.Range("N51").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O51").Value = "Provincia de Diguillín"
$ws.Range("P51").Value = 575
$ws.Range("Q51").Value = 1
$ws.Range("R51").Value = "Hortaliza"
